$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.604.70"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.918.61"
$ws.Range("E3").Value = "  +1.29%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "469.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.96%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.40%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.02%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.746"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.88%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.11%  "

# Row 11 - ShibaInu
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000318"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.17%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.85%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.558.04"
$ws.Range("E13").Value = "  +1.93%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15 - Uniswap
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.48%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.960.27"
$ws.Range("E16").Value = "  +2.71%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.29%  "

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +7.48%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "68.068.86"
$ws.Range("E20").Value = "  +1.30%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.67%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.12%  "

# Row 23 - ImmutableX
$ws.Range("E23").Value = "  +8.40%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.31%  "

# Row 25 - PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.22%  "

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "

# Row 27 - RenderToken
$ws.Range("E27").Value = "  +13.87%  "

# Row 28 - Filecoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.67%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +2.29%  "

# Row 30 - Bittensor
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "736.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.135"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.23%  "

# Row 32 - Cosmos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "

# Row 33 - Toncoin
$ws.Range("E33").Value = "  -0.60%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.32%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +7.20%  "

# Row 36 - OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.11%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +3.81%  "

# Row 38 - Dai
$ws.Range("E38").Value = "  +0.15%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +4.34%  "

# Row 40 - ThetaToken
$ws.Range("E40").Value = "  +1.79%  "

# Row 41 - PEPE
$ws.Range("D41").Value = "0.0₃0694"
$ws.Range("E41").Value = "  -8.67%  "

# Row 42 - TheGraph
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.347"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.09%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  +4.94%  "

# Row 44 - LidoDAOToken
$ws.Range("E44").Value = "  +4.33%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.03%  "

# Row 46 - Fetch.AI
$ws.Range("E46").Value = "  +13.82%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  +8.09%  "

# Row 48 - ARBITRUM -> ApeXProtocol
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "

# Row 49 - ApeXProtocol -> ARBITRUM
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.20%  "

# Row 50 - Stacks -> Monero
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.41%  "

# Row 51 - Monero -> Stacks
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.59%  "
